$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Split the run "left-bottom" into "left-" (keeps the original
#    run formatting) + "middle" (a new run with minimal / default
#    run properties, matching the target XML).
# ------------------------------------------------------------------

# Narrow down to the "bottom" substring inside "left-bottom" and
# delete it, leaving "left-" behind with its original run formatting
# untouched.
$r1 = $d.Content
$found1 = $r1.Find.Execute("left-bottom", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Start = $r1.End - 6
$r1.Delete()

# Grab a "template" range elsewhere in the document whose run
# properties are already minimal (just rtl) -- the single-letter "E"
# run inside "Extraction" -- and copy its FormattedText (this carries
# the run-level formatting only, without touching the source text).
$full = $d.Content.Text
$tmplIdx = $full.IndexOf("xtraction")
$srcChar = $d.Range($tmplIdx - 1, $tmplIdx)
$ft = $srcChar.FormattedText

# Re-find "left-" (now that "bottom" is gone) and collapse the range
# to its end, then drop in the minimally-formatted template text
# right after it.
$r2 = $d.Content
$found2 = $r2.Find.Execute("left-", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Collapse(0)
$r2.FormattedText = $ft

# The pasted text is still "E" at this point (the template's text) --
# retarget just that freshly-inserted character and change it to
# "middle" without touching the original template run.
$full2 = $d.Content.Text
$newIdx = $full2.IndexOf("left-E")
$target = $d.Range($newIdx + 5, $newIdx + 6)
$target.Text = "middle"

# ------------------------------------------------------------------
# 2) Add a footer distance of 720 twips (36pt) to the section's page
#    margins.
# ------------------------------------------------------------------
$d.Sections.Item(1).PageSetup.FooterDistance = 36
